$wb = $excel.ActiveWorkbook

$wsRanking = $wb.Worksheets.Item("Ranking")
$wsPvalores = $wb.Worksheets.Item("Matriz_Pvalores")
$wsDM = $wb.Worksheets.Item("Matriz_DM_Original")

# --- Sheet "Ranking": CRPS_Mean / CRPS_Median updates for LSPMW (row 5) and EnCQR-LSTM (row 9) ---
$wsRanking.Range("H5").Value = [double]"0.004388101543955188"
$wsRanking.Range("I5").Value = [double]"0.002535966299040615"
$wsRanking.Range("H9").Value = [double]"0.06657783537960374"
$wsRanking.Range("I9").Value = [double]"0.06763749094517246"

# --- Sheet "Matriz_Pvalores": recomputed p-values touching LSPMW (row/col E,5) and EnCQR-LSTM (row/col J,10) ---
$wsPvalores.Range("E2").Value = [double]"1.510615714961361e-06"
$wsPvalores.Range("J2").Value = [double]"0.006872943992138003"

$wsPvalores.Range("E3").Value = [double]"0.09586215660673125"
$wsPvalores.Range("J3").Value = [double]"3.01254876822199e-07"

$wsPvalores.Range("E4").Value = [double]"0.02012931258132666"
$wsPvalores.Range("J4").Value = [double]"3.217468258487344e-07"

$wsPvalores.Range("B5").Value = [double]"1.510615714961361e-06"
$wsPvalores.Range("C5").Value = [double]"0.09586215660673125"
$wsPvalores.Range("D5").Value = [double]"0.02012931258132666"
$wsPvalores.Range("F5").Value = [double]"1.803747326389882e-06"
$wsPvalores.Range("G5").Value = [double]"0.01288173122514502"
$wsPvalores.Range("H5").Value = [double]"0.02457370594547625"
$wsPvalores.Range("I5").Value = [double]"0.3137386149345722"
$wsPvalores.Range("J5").Value = [double]"3.218905919588622e-07"

$wsPvalores.Range("E6").Value = [double]"1.803747326389882e-06"
$wsPvalores.Range("J6").Value = [double]"0.05325072225270833"

$wsPvalores.Range("E7").Value = [double]"0.01288173122514502"
$wsPvalores.Range("J7").Value = [double]"2.629998243697251e-05"

$wsPvalores.Range("E8").Value = [double]"0.02457370594547625"
$wsPvalores.Range("J8").Value = [double]"2.731430908720078e-05"

$wsPvalores.Range("E9").Value = [double]"0.3137386149345722"
$wsPvalores.Range("J9").Value = [double]"2.179381251021795e-07"

$wsPvalores.Range("B10").Value = [double]"0.006872943992138003"
$wsPvalores.Range("C10").Value = [double]"3.01254876822199e-07"
$wsPvalores.Range("D10").Value = [double]"3.217468258487344e-07"
$wsPvalores.Range("E10").Value = [double]"3.218905919588622e-07"
$wsPvalores.Range("F10").Value = [double]"0.05325072225270833"
$wsPvalores.Range("G10").Value = [double]"2.629998243697251e-05"
$wsPvalores.Range("H10").Value = [double]"2.731430908720078e-05"
$wsPvalores.Range("I10").Value = [double]"2.179381251021795e-07"

# --- Sheet "Matriz_DM_Original": recomputed DM statistics touching LSPMW (row/col E,5) and EnCQR-LSTM (row/col J,10) ---
$wsDM.Range("E2").Value = [double]"18.69782472683895"
$wsDM.Range("J2").Value = [double]"4.030940872175816"

$wsDM.Range("E3").Value = [double]"-1.973663698962231"
$wsDM.Range("J3").Value = [double]"-24.5382785568941"

$wsDM.Range("E4").Value = [double]"-3.137595145501581"
$wsDM.Range("J4").Value = [double]"-24.26826760023854"

$wsDM.Range("B5").Value = [double]"-18.69782472683895"
$wsDM.Range("C5").Value = [double]"1.973663698962231"
$wsDM.Range("D5").Value = [double]"3.137595145501581"
$wsDM.Range("F5").Value = [double]"-18.1449514972347"
$wsDM.Range("G5").Value = [double]"-3.496456881586686"
$wsDM.Range("H5").Value = [double]"-2.981977649605745"
$wsDM.Range("I5").Value = [double]"1.09939767802291"
$wsDM.Range("J5").Value = [double]"-24.26644468315485"

$wsDM.Range("E6").Value = [double]"18.1449514972347"
$wsDM.Range("J6").Value = [double]"-2.400608331026901"

$wsDM.Range("E7").Value = [double]"3.496456881586686"
$wsDM.Range("J7").Value = [double]"-11.47502573868122"

$wsDM.Range("E8").Value = [double]"2.981977649605745"
$wsDM.Range("J8").Value = [double]"-11.40000679487877"

$wsDM.Range("E9").Value = [double]"-1.09939767802291"
$wsDM.Range("J9").Value = [double]"-25.91021298372895"

$wsDM.Range("B10").Value = [double]"-4.030940872175816"
$wsDM.Range("C10").Value = [double]"24.5382785568941"
$wsDM.Range("D10").Value = [double]"24.26826760023854"
$wsDM.Range("E10").Value = [double]"24.26644468315485"
$wsDM.Range("F10").Value = [double]"2.400608331026901"
$wsDM.Range("G10").Value = [double]"11.47502573868122"
$wsDM.Range("H10").Value = [double]"11.40000679487877"
$wsDM.Range("I10").Value = [double]"25.91021298372895"
